$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-04 Saturday" "2025-10-05 Sunday"

Replace-Text "498×6=" "209×4="
Replace-Text "600×9=" "523×8="
Replace-Text "858×7=" "737×2="
Replace-Text "655×3=" "744×2="
Replace-Text "440×2=" "270×2="
Replace-Text "430×6=" "571×7="
Replace-Text "754×5=" "976×7="
Replace-Text "433×6=" "132×3="
Replace-Text "584×7=" "245×9="
Replace-Text "951×4=" "540×3="
Replace-Text "799×8=" "695×8="
Replace-Text "687×5=" "520×9="
Replace-Text "637×4=" "565×8="
Replace-Text "859×3=" "202×8="
Replace-Text "763×5=" "440×4="
Replace-Text "920×2=" "352×8="
Replace-Text "647×9=" "336×2="
Replace-Text "889×7=" "310×4="
Replace-Text "443×9=" "135×5="
Replace-Text "292×9=" "266×8="
Replace-Text "743×2=" "681×9="
Replace-Text "365×2=" "996×2="
Replace-Text "426×5=" "830×5="
Replace-Text "489×4=" "238×6="
Replace-Text "815×6=" "479×5="
